# RTE_Configurator test plan - add new test case rows for
# "All events of type <OPERATION-INVOKED-EVENT> ... must not be mapped"
# on the "Overview" worksheet (sheet2), inserted right before the
# "RTECONFIG.1 / RTECONFIG.2" rows (originally rows 41-42, now 43-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert two new blank rows at 41-42; everything below (old rows 41-43)
# shifts down to 43-45, and all ranges that reference those rows (merged
# cells, dimension, autofilter-adjacent ranges, etc.) shift accordingly.
$ws.Rows("41:42").Insert()

# ---- Column A: requirement id "RTE.EXCEPT.OIE" (merged A41:A42) ----
$ws.Range("A41").Value = "RTE.EXCEPT.OIE"
$ws.Range("A18").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A41").WrapText = $false

$ws.Range("A24").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A42").WrapText = $false

$ws.Range("A41:A42").Merge()

# ---- Column B: requirement description (merged B41:B42) ----
$ws.Range("B41").Value = "All the events of type <OPERATION-INVOKED-EVENT> which belongs to ASWC IoHwAb  or ASWC Aswc_IntDcm mult not be mapped"
$ws.Range("B18").Copy()
$ws.Range("B41").PasteSpecial(-4122)

$ws.Range("B24").Copy()
$ws.Range("B42").PasteSpecial(-4122)

$ws.Range("B41:B42").Merge()

# ---- Column C: the two test descriptions (rich text, "Test 1:" bold) ----
$ws.Range("C41").Value = "Test 1: Provide to the tool one file havin one event of type <OPERATION-INVOKED-EVENT> assigned to ASWC IoHwAb"
$ws.Range("C42").Value = "Test 1: Provide to the tool one file havin one event of type <OPERATION-INVOKED-EVENT> assigned to ASWC Aswc_IntDcm"

$ws.Range("C40").Copy()
$ws.Range("C41").PasteSpecial(-4122)
$ws.Range("C42").PasteSpecial(-4122)

$ws.Range("C41").Characters(1, 7).Font.Bold = $true
$ws.Range("C42").Characters(1, 7).Font.Bold = $true

# ---- Column D: test result description ----
$ws.Range("D41").Value = "- check that in the output file the event is not found"
$ws.Range("D42").Value = "- check that in the output file the event is not found"

$ws.Range("D40").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("D42").PasteSpecial(-4122)

# ---- Column E: Implemented = done ----
$ws.Range("E41").Value = "done"
$ws.Range("E42").Value = "done"

$ws.Range("E40").Copy()
$ws.Range("E41").PasteSpecial(-4122)
$ws.Range("E42").PasteSpecial(-4122)

# ---- Column F: Passed = YES ----
$ws.Range("F41").Value = "YES"
$ws.Range("F42").Value = "YES"

$ws.Range("F40").Copy()
$ws.Range("F41").PasteSpecial(-4122)
$ws.Range("F42").PasteSpecial(-4122)

# ---- Row heights ----
$ws.Rows.Item(41).RowHeight = 45
$ws.Rows.Item(42).RowHeight = 30

# ---- Re-extend conditional formatting / data validation ranges that
# should now cover the two extra rows (Excel normally grows these
# automatically when rows are inserted inside their range; make sure
# they end up covering through row 44/45 as expected). ----
$ws.Range("E2:E44").FormatConditions.Delete()
$ws.Range("E2:E44").FormatConditions.Add(8, 3, """done""")
$ws.Range("E2:E44").FormatConditions.Item(1).Interior.Color = 13561798
$ws.Range("E2:E44").FormatConditions.Item(1).Font.Color = 6299648
$ws.Range("E2:E44").FormatConditions.Add(8, 3, """not done""")
$ws.Range("E2:E44").FormatConditions.Item(2).Interior.Color = 13551615
$ws.Range("E2:E44").FormatConditions.Item(2).Font.Color = 402

$ws.Range("E46:E1048576").FormatConditions.Delete()
$ws.Range("E46:E1048576").FormatConditions.Add(8, 3, """done""")
$ws.Range("E46:E1048576").FormatConditions.Item(1).Interior.Color = 13561798
$ws.Range("E46:E1048576").FormatConditions.Item(1).Font.Color = 6299648
$ws.Range("E46:E1048576").FormatConditions.Add(8, 3, """not done""")
$ws.Range("E46:E1048576").FormatConditions.Item(2).Interior.Color = 13551615
$ws.Range("E46:E1048576").FormatConditions.Item(2).Font.Color = 402

$ws.Range("F2:F44").FormatConditions.Delete()
$ws.Range("F2:F44").FormatConditions.Add(8, 3, """N/A""")
$ws.Range("F2:F44").FormatConditions.Item(1).Interior.Color = 10284031
$ws.Range("F2:F44").FormatConditions.Item(1).Font.Color = 25600
$ws.Range("F2:F44").FormatConditions.Add(8, 3, """NO""")
$ws.Range("F2:F44").FormatConditions.Item(2).Interior.Color = 13551615
$ws.Range("F2:F44").FormatConditions.Item(2).Font.Color = 402
$ws.Range("F2:F44").FormatConditions.Add(8, 3, """YES""")
$ws.Range("F2:F44").FormatConditions.Item(3).Interior.Color = 13561798
$ws.Range("F2:F44").FormatConditions.Item(3).Font.Color = 6299648

$ws.Range("E2:E44").Validation.Delete()
$ws.Range("E2:E44").Validation.Add(3, 1, 1, "done, not done")
$ws.Range("F2:F44").Validation.Delete()
$ws.Range("F2:F44").Validation.Add(3, 1, 1, "YES, NO, N/A")
